$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.607.52'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '2.612.03'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '538.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.75%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.54'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("E10").Value = '  +1.22%  '
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("D13").Value = '3.067.82'
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("D14").Value = '59.531.41'
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.630.06'
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '341.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.95%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("E27").Value = '  +2.69%  '
$ws.Range("E28").Value = '  +3.00%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  +5.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.79'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '150.77'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("E35").Value = '  +0.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.844'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.53%  '
$ws.Range("E37").Value = '  -0.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.828'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '278.13'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0524'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.86%  '
$ws.Range("D46").Value = '1.941.63'
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.48'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.14%  '
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = '  +1.76%  '
